$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 4.127494786624378
$ws.Range("E2").Value = 0.40788
$ws.Range("F2").Value = 3.719614786624378

$ws.Range("D3").Value = 3.955052113479161
$ws.Range("E3").Value = 0.66693
$ws.Range("F3").Value = 3.288122113479161

$ws.Range("D4").Value = 3.805243826469577
$ws.Range("E4").Value = 1.01763
$ws.Range("F4").Value = 2.787613826469577

$ws.Range("D5").Value = 3.777562248922468
$ws.Range("E5").Value = 1.30713
$ws.Range("F5").Value = 2.470432248922468

$ws.Range("D6").Value = 3.761102860317084
$ws.Range("E6").Value = 1.46323
$ws.Range("F6").Value = 2.297872860317085

$ws.Range("D7").Value = 3.732473588889856
$ws.Range("E7").Value = 1.60923
$ws.Range("F7").Value = 2.123243588889856

$ws.Range("D8").Value = 3.707770631891534
$ws.Range("E8").Value = 1.74012
$ws.Range("F8").Value = 1.967650631891535

$ws.Range("D9").Value = 3.689215858993487
$ws.Range("E9").Value = 1.86612
$ws.Range("F9").Value = 1.823095858993487

$ws.Range("D10").Value = 3.668793200867766
$ws.Range("E10").Value = 1.93652
$ws.Range("F10").Value = 1.732273200867766
